$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("B2").Value = 0.5
    $ws.Range("C2").Value = 0.5384615384615384
    $ws.Range("D2").Value = 0.5185185185185186
    $ws.Range("E2").Value = 13
    $ws.Range("B3").Value = 0.4
    $ws.Range("C3").Value = 0.3636363636363636
    $ws.Range("D3").Value = 0.380952380952381
    $ws.Range("E3").Value = 11
    $ws.Range("B5").Value = 0.45
    $ws.Range("C5").Value = 0.451048951048951
    $ws.Range("D5").Value = 0.4497354497354498
    $ws.Range("B6").Value = 0.4541666666666667
    $ws.Range("D6").Value = 0.4554673721340389
    $ws.Range("B7").Value = 0.4444444444444444
    $ws.Range("C7").Value = 0.6153846153846154
    $ws.Range("D7").Value = 0.5161290322580646
    $ws.Range("E7").Value = 13
    $ws.Range("B8").Value = 0.1666666666666667
    $ws.Range("C8").Value = 0.09090909090909091
    $ws.Range("D8").Value = 0.1176470588235294
    $ws.Range("E8").Value = 11
    $ws.Range("B9").Value = 0.375
    $ws.Range("C9").Value = 0.375
    $ws.Range("D9").Value = 0.375
    $ws.Range("E9").Value = 0.375
    $ws.Range("B10").Value = 0.3055555555555555
    $ws.Range("C10").Value = 0.3531468531468532
    $ws.Range("D10").Value = 0.316888045540797
    $ws.Range("B11").Value = 0.3171296296296296
    $ws.Range("C11").Value = 0.375
    $ws.Range("D11").Value = 0.3334914611005693
    $ws.Range("B12").Value = 0.5454545454545454
    $ws.Range("C12").Value = 0.4615384615384616
    $ws.Range("D12").Value = 0.4999999999999999
    $ws.Range("E12").Value = 13
    $ws.Range("B13").Value = 0.4615384615384616
    $ws.Range("C13").Value = 0.5454545454545454
    $ws.Range("D13").Value = 0.4999999999999999
    $ws.Range("E13").Value = 11
    $ws.Range("B15").Value = 0.5034965034965035
    $ws.Range("C15").Value = 0.5034965034965035
    $ws.Range("D15").Value = 0.4999999999999999
    $ws.Range("B16").Value = 0.506993006993007
    $ws.Range("D16").Value = 0.4999999999999998
    $ws.Range("B17").Value = 0.6153846153846154
    $ws.Range("C17").Value = 0.6153846153846154
    $ws.Range("D17").Value = 0.6153846153846154
    $ws.Range("E17").Value = 13
    $ws.Range("B18").Value = 0.5454545454545454
    $ws.Range("C18").Value = 0.5454545454545454
    $ws.Range("D18").Value = 0.5454545454545454
    $ws.Range("E18").Value = 11
    $ws.Range("B20").Value = 0.5804195804195804
    $ws.Range("C20").Value = 0.5804195804195804
    $ws.Range("D20").Value = 0.5804195804195804
    $ws.Range("B21").Value = 0.5833333333333334
    $ws.Range("D21").Value = 0.5833333333333334
    $ws.Range("B22").Value = 0.5384615384615384
    $ws.Range("C22").Value = 0.5384615384615384
    $ws.Range("D22").Value = 0.5384615384615384
    $ws.Range("E22").Value = 13
    $ws.Range("B23").Value = 0.4545454545454545
    $ws.Range("C23").Value = 0.4545454545454545
    $ws.Range("D23").Value = 0.4545454545454545
    $ws.Range("E23").Value = 11
    $ws.Range("B24").Value = 0.5
    $ws.Range("C24").Value = 0.5
    $ws.Range("D24").Value = 0.5
    $ws.Range("E24").Value = 0.5
    $ws.Range("B25").Value = 0.4965034965034965
    $ws.Range("C25").Value = 0.4965034965034965
    $ws.Range("D25").Value = 0.4965034965034965
    $ws.Range("B26").Value = 0.5
    $ws.Range("C26").Value = 0.5
    $ws.Range("D26").Value = 0.5
